$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the existing header cell (H1) onto the new header
# cells I1:J1 so they pick up the same bold/border/centered style (s="1"),
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in new data columns I (I0) and J (IF) for rows 2-12
$data = @(
    @(2, 2, 3),
    @(3, 6, 6),
    @(4, 6, 7),
    @(5, 9, 9),
    @(6, 8, 8),
    @(7, 10, 10),
    @(8, 5, 7),
    @(9, 5, 6),
    @(10, 5, 7),
    @(11, 7, 7),
    @(12, 1, 2)
)

foreach ($row in $data) {
    $r = $row[0]
    $iVal = $row[1]
    $jVal = $row[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
